$wb = $excel.ActiveWorkbook

# --- "Page" sheet (sheet3.xml): insert two new lead-in rows above the
#     existing data table, shifting the rest of the table down by two rows.
$ws = $wb.Worksheets.Item("Page")

$ws.Range("2:3").Insert()

# New row 2: a "Paypyrus" splash/title row
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = "Paypyrus"
$ws.Cells.Item(2, 6).Value = "Tap anywhere to continue!"

# New row 3: a "Get started with the basics" intro row
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = "Get started with the basics"
$ws.Cells.Item(3, 6).Value = "Start out right ( to left!!!) with these foundational lessons to get you going in the right (left) direction!"

# Match the author's final cursor position on this sheet.
$ws.Range("F12").Select()

# The "Page" tab (third tab) ends up the active tab in the saved workbook.
$ws.Activate()
